$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$changes = @{
    4  = -6
    5  = -1
    10 = -2
    11 = -7
    12 = -2
    15 = -7
    17 = 5
    23 = -5
    24 = 2
    27 = 1
    29 = 7
    30 = -7
    33 = -8
    35 = -2
    37 = -5
    40 = -7
    42 = -1
    48 = -6
    53 = -2
    57 = -1
    58 = -2
    59 = -2
    62 = -2
    63 = -5
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
